$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, then a hashtable of Column -> new text value.
# Source of truth: the cell-level unified diff for xl/worksheets/sheet1.xml.
$updates = @(
    @{ Row=2; Cells=@{ "D"="278.65"; "E"="4.92%"; "G"="19" } }
    @{ Row=3; Cells=@{ "D"="26.94"; "E"="0.70%"; "G"="19" } }
    @{ Row=4; Cells=@{ "D"="4.934"; "E"="5.10%"; "G"="19" } }
    @{ Row=5; Cells=@{ "D"="0.06405"; "E"="5.31%"; "G"="19" } }
    @{ Row=6; Cells=@{ "D"="7.002"; "E"="4.58%"; "G"="19" } }
    @{ Row=7; Cells=@{ "D"="3.360"; "E"="5.97%"; "G"="19" } }
    @{ Row=8; Cells=@{ "D"="0.8886"; "E"="4.53%"; "G"="19" } }
    @{ Row=9; Cells=@{ "D"="1.050"; "E"="16.20%"; "G"="19" } }
    @{ Row=10; Cells=@{ "E"="5.96%"; "G"="19" } }
    @{ Row=11; Cells=@{ "D"="0.05240"; "E"="8.96%"; "G"="19" } }
    @{ Row=12; Cells=@{ "D"="0.07403"; "E"="4.28%"; "G"="19" } }
    @{ Row=13; Cells=@{ "D"="0.03100"; "E"="-2.04%"; "G"="19" } }
    @{ Row=14; Cells=@{ "D"="0.09068"; "E"="0.55%"; "G"="19" } }
    @{ Row=15; Cells=@{ "D"="0.001578"; "E"="2.81%"; "G"="19" } }
    @{ Row=16; Cells=@{ "D"="0.0006372"; "E"="4.98%"; "G"="19" } }
    @{ Row=17; Cells=@{ "D"="0.006033"; "E"="0.60%"; "G"="19" } }
    @{ Row=18; Cells=@{ "D"="3.496"; "E"="1.12%"; "G"="19" } }
    @{ Row=19; Cells=@{ "D"="2.295"; "E"="0.74%"; "G"="19" } }
    @{ Row=20; Cells=@{ "G"="19" } }
    @{ Row=21; Cells=@{ "D"="0.1334"; "E"="2.60%"; "G"="19" } }
    @{ Row=22; Cells=@{ "D"="3.911"; "E"="-4.18%"; "G"="19" } }
    @{ Row=23; Cells=@{ "D"="0.04357"; "E"="2.86%"; "G"="19" } }
    @{ Row=24; Cells=@{ "D"="0.001183"; "E"="0.18%"; "G"="19" } }
    @{ Row=25; Cells=@{ "D"="0.003688"; "E"="-10.70%"; "G"="19" } }
    @{ Row=26; Cells=@{ "D"="0.0001203"; "E"="0.19%"; "G"="19" } }
    @{ Row=27; Cells=@{ "E"="1.01%"; "G"="19" } }
    @{ Row=28; Cells=@{ "G"="19" } }
    @{ Row=29; Cells=@{ "G"="19" } }
    @{ Row=30; Cells=@{ "G"="19" } }
    @{ Row=31; Cells=@{ "G"="19" } }
    @{ Row=32; Cells=@{ "G"="19" } }
    @{ Row=33; Cells=@{ "G"="19" } }
    @{ Row=34; Cells=@{ "G"="19" } }
    @{ Row=35; Cells=@{ "G"="19" } }
    @{ Row=36; Cells=@{ "G"="19" } }
    @{ Row=37; Cells=@{ "G"="19" } }
    @{ Row=38; Cells=@{ "G"="19" } }
    @{ Row=39; Cells=@{ "G"="19" } }
    @{ Row=40; Cells=@{ "D"="0.04077"; "E"="4.15%"; "G"="19" } }
    @{ Row=41; Cells=@{ "D"="0.006653"; "E"="58.40%"; "G"="19" } }
    @{ Row=42; Cells=@{ "D"="0.1177"; "E"="5.54%"; "G"="19" } }
    @{ Row=43; Cells=@{ "E"="12.07%"; "G"="19" } }
    @{ Row=44; Cells=@{ "D"="0.01290"; "E"="12.30%"; "G"="19" } }
    @{ Row=45; Cells=@{ "D"="0.00005279"; "E"="2.94%"; "G"="19" } }
    @{ Row=46; Cells=@{ "E"="-0.02%"; "G"="19" } }
    @{ Row=47; Cells=@{ "E"="1,612.21%"; "G"="19" } }
    @{ Row=48; Cells=@{ "D"="0.02125"; "E"="-13.20%"; "G"="19" } }
    @{ Row=49; Cells=@{ "E"="-0.02%"; "G"="19" } }
    @{ Row=50; Cells=@{ "E"="-0.09%"; "G"="19" } }
    @{ Row=51; Cells=@{ "G"="19" } }
)

foreach ($update in $updates) {
    foreach ($col in $update.Cells.Keys) {
        $addr = "$col$($update.Row)"
        $cell = $ws.Range($addr)
        # Force text storage so e.g. "19" / "4.92%" are not reinterpreted
        # as a number/percentage by Excel's input parser.
        $cell.NumberFormat = "@"
        $cell.Value = $update.Cells[$col]
        # Drop back to the default style so no stray format/style
        # attribute is left behind on the cell (matches source diff,
        # which only touches the <t> text, not any style ids).
        $cell.Style = "Normal"
    }
}
